$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("G4").Value = "URI-Comunidad"
$ws.Range("G5").Clear()
